$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert 5 new rows (13-17) to make room for the new
# "Docentes responsaveis" block and shift the rest of the table down.
$ws.Rows("13:17").Insert()

# Step 2: Column A now gets its own column definition instead of being
# merged with column B (touching column B forces the engine to split the
# shared col-range while keeping column A at its original exact width).
$ws.Columns(2).ColumnWidth = $ws.Columns(2).ColumnWidth

# Step 3: Write every cell value so the whole table matches the new text.
$ws.Cells.Item(1,2).Value2 = "Ementa atual:"
$ws.Cells.Item(1,3).Value2 = "Ementa modificada (dados modificados em vermelho):"
$ws.Cells.Item(2,2).Value2 = "LOM3081"
$ws.Cells.Item(2,3).Value2 = "LOM3081"
$ws.Cells.Item(3,1).Value2 = "Nome:"
$ws.Cells.Item(3,2).Value2 = " Introdução à Mecânica dos Sólidos"
$ws.Cells.Item(3,3).Value2 = " Introdução à Mecânica dos Sólidos"
$ws.Cells.Item(4,1).Value2 = "Name:"
$ws.Cells.Item(4,2).Value2 = "Introduction to Solid Mechanics"
$ws.Cells.Item(4,3).Value2 = "Introduction to Solid Mechanics"
$ws.Cells.Item(5,1).Value2 = "Créditos-aula:"
$ws.Cells.Item(5,2).Value2 = "2"
$ws.Cells.Item(5,3).Value2 = "2"
$ws.Cells.Item(6,1).Value2 = "Créditos-trabalho"
$ws.Cells.Item(6,2).Value2 = "0"
$ws.Cells.Item(6,3).Value2 = "0"
$ws.Cells.Item(7,1).Value2 = "Carga horária:"
$ws.Cells.Item(7,2).Value2 = "30 h"
$ws.Cells.Item(7,3).Value2 = "30 h"
$ws.Cells.Item(8,1).Value2 = "Ativação:"
$ws.Cells.Item(8,2).Value2 = "01/01/2012"
$ws.Cells.Item(8,3).Value2 = "01/01/2012"
$ws.Cells.Item(9,1).Value2 = "Semestre ideal:"
$ws.Cells.Item(9,2).Value2 = "EF-4,EA-4,EP-6,EQD-4,EQN-6"
$ws.Cells.Item(9,3).Value2 = "EF-4,EA-4,EP-6,EQD-4,EQN-6"
$ws.Cells.Item(10,1).Value2 = "Objetivos:"
$ws.Cells.Item(10,2).Value2 = "Fornecer conceitos relacionados ao comportamento dos sólidos deformáveis, capacitando ao cálculo de tensões e deformações em sistemas de barras axialmente carregadas, à análise dos estados planos de tensão e deformação, bem como prover o conhecimento e a aplicação das propriedades elásticas dos materiais."
$ws.Cells.Item(10,3).Value2 = "Fornecer conceitos relacionados ao comportamento dos sólidos deformáveis, capacitando ao cálculo de tensões e deformações em sistemas de barras axialmente carregadas, à análise dos estados planos de tensão e deformação, bem como prover o conhecimento e a aplicação das propriedades elásticas dos materiais."
$ws.Cells.Item(11,1).Value2 = "Objectives:"
$ws.Cells.Item(12,1).Value2 = "Docentes responsáveis:"
$ws.Cells.Item(13,2).Value2 = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Cells.Item(13,3).Value2 = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Cells.Item(14,2).Value2 = "3480026 - João Paulo Pascon"
$ws.Cells.Item(14,3).Value2 = "3480026 - João Paulo Pascon"
$ws.Cells.Item(15,2).Value2 = "5840793 - Sérgio Schneider"
$ws.Cells.Item(15,3).Value2 = "5840793 - Sérgio Schneider"
$ws.Cells.Item(16,2).Value2 = "7797767 - Viktor Pastoukhov"
$ws.Cells.Item(16,3).Value2 = "7797767 - Viktor Pastoukhov"
$ws.Cells.Item(17,1).Value2 = "Programa resumido:"
$ws.Cells.Item(17,2).Value2 = "Considerações fundamentais; Tensão e deformação em membros carregados axialmente; Análise de tensão e deformação; Relações tensão-deformação no regime elástico."
$ws.Cells.Item(17,3).Value2 = "Considerações fundamentais; Tensão e deformação em membros carregados axialmente; Análise de tensão e deformação; Relações tensão-deformação no regime elástico."
$ws.Cells.Item(18,1).Value2 = "Short syllabus:"
$ws.Cells.Item(19,1).Value2 = "Programa:"
$ws.Cells.Item(19,2).Value2 = "1.Considerações Fundamentais: Propósito da Mecânica dos Sólidos; Carregamentos e Esforços Solicitantes; Tensão Normal e Tensão Cisalhante; Tensões admissíveis.`n2.Tensão e Deformação em Membros Carregados Axialmente: Elasticidade linear e o Módulo de Young, Sistemas Isostáticos e Hiperestáticos; Efeitos da Temperatura.`n3.Análise de Tensão e Deformação: Variação da Tensão com o Plano de Corte; Estado Plano de Tensão; Tensões Principais e Máxima Tensão de Cisalhamento; O Círculo de Mohr para Tensão Plana; Tensão Triaxial; Deformação Angular e Módulo de Elasticidade Transversal; Coeficiente de Poisson; Transformação do Estado Plano de Deformação.`n4.Relações Tensão-Deformação no Regime Elástico: Elasticidade, Homogeneidade e Isotropia; Lei de Hooke para Tensão Triaxial em Materiais Isotrópicos; Relações entre as Constantes Elásticas; Aplicação em Vasos de Pressão de Paredes Finas."
$ws.Cells.Item(19,3).Value2 = "1.Considerações Fundamentais: Propósito da Mecânica dos Sólidos; Carregamentos e Esforços Solicitantes; Tensão Normal e Tensão Cisalhante; Tensões admissíveis.`n2.Tensão e Deformação em Membros Carregados Axialmente: Elasticidade linear e o Módulo de Young, Sistemas Isostáticos e Hiperestáticos; Efeitos da Temperatura.`n3.Análise de Tensão e Deformação: Variação da Tensão com o Plano de Corte; Estado Plano de Tensão; Tensões Principais e Máxima Tensão de Cisalhamento; O Círculo de Mohr para Tensão Plana; Tensão Triaxial; Deformação Angular e Módulo de Elasticidade Transversal; Coeficiente de Poisson; Transformação do Estado Plano de Deformação.`n4.Relações Tensão-Deformação no Regime Elástico: Elasticidade, Homogeneidade e Isotropia; Lei de Hooke para Tensão Triaxial em Materiais Isotrópicos; Relações entre as Constantes Elásticas; Aplicação em Vasos de Pressão de Paredes Finas."
$ws.Cells.Item(20,1).Value2 = "Syllabus:"
$ws.Cells.Item(21,1).Value2 = "Avaliação:"
$ws.Cells.Item(22,1).Value2 = "Método:"
$ws.Cells.Item(22,2).Value2 = "Para compor a Nota no Semestre (NS) serão feitas duas avaliações (P1 e P2) e o critério de cálculo será: NS = (P1 + P2)/2."
$ws.Cells.Item(22,3).Value2 = "Para compor a Nota no Semestre (NS) serão feitas duas avaliações (P1 e P2) e o critério de cálculo será: NS = (P1 + P2)/2."
$ws.Cells.Item(23,1).Value2 = "Critério:"
$ws.Cells.Item(23,2).Value2 = "Serão considerados aprovados os alunos que obtiverem NS maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem NS menor que 3,0. Para os alunos que obtiverem NS maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$ws.Cells.Item(23,3).Value2 = "Serão considerados aprovados os alunos que obtiverem NS maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem NS menor que 3,0. Para os alunos que obtiverem NS maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$ws.Cells.Item(24,1).Value2 = "Norma de recuperação:"
$ws.Cells.Item(24,2).Value2 = "A prova de Recuperação (R) irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
$ws.Cells.Item(24,3).Value2 = "A prova de Recuperação (R) irá compor a nota final (NF) da seguinte forma: NF = (R + NS)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
$ws.Cells.Item(25,1).Value2 = "Bibliografia:"
$ws.Cells.Item(25,2).Value2 = "1. J.M. GERE. Mecânica dos Materiais. São Paulo: Pioneira Thomson Learning, 2003, 698p.`n2. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF. Resistência dos Materiais. São Paulo: McGraw Hill. 4a Ed., 2006, 758p.`n3. R.R. CRAIG,Jr. Mecânica dos Materiais. Rio de Janeiro LTC. 2a Ed., 2003, 552p. `n4. R.C. HIBBELER. Resistência dos Materiais. São Paulo: Pearson Prentice Hall. 5a Ed., 2006, 670p.`n5. A.C. UGURAL. Mecânica dos Materiais. Rio de Janeiro LTC, 2009, 638p.`n6. A.R. RAGAB, S.E. BAYOUMI. Engineering Solid Mechanics, Fundamentals and Applications. New York: CRC Press, 1999, 921p. `n7. POPOV, E. P. Introdução à Mecânica dos Sólidos, São Paulo: Edgard Blücher, 1978, 552p.`n8. A. HIGDON, E.H. OHLSEN, W.B. STILES, J.A. WEESE, W.F. RILEY. Mecânica dos Materiais.  Rio de Janeiro: Guanabara Dois. 3a Ed., 1981, 549p."
$ws.Cells.Item(25,3).Value2 = "1. J.M. GERE. Mecânica dos Materiais. São Paulo: Pioneira Thomson Learning, 2003, 698p.`n2. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF. Resistência dos Materiais. São Paulo: McGraw Hill. 4a Ed., 2006, 758p.`n3. R.R. CRAIG,Jr. Mecânica dos Materiais. Rio de Janeiro LTC. 2a Ed., 2003, 552p. `n4. R.C. HIBBELER. Resistência dos Materiais. São Paulo: Pearson Prentice Hall. 5a Ed., 2006, 670p.`n5. A.C. UGURAL. Mecânica dos Materiais. Rio de Janeiro LTC, 2009, 638p.`n6. A.R. RAGAB, S.E. BAYOUMI. Engineering Solid Mechanics, Fundamentals and Applications. New York: CRC Press, 1999, 921p. `n7. POPOV, E. P. Introdução à Mecânica dos Sólidos, São Paulo: Edgard Blücher, 1978, 552p.`n8. A. HIGDON, E.H. OHLSEN, W.B. STILES, J.A. WEESE, W.F. RILEY. Mecânica dos Materiais.  Rio de Janeiro: Guanabara Dois. 3a Ed., 1981, 549p."
$ws.Cells.Item(26,1).Value2 = "Requisitos:"
$ws.Cells.Item(27,2).Value2 = "LOM3257 -  Mecânica Clássica  (Requisito fraco)`n"
$ws.Cells.Item(27,3).Value2 = "LOM3257 -  Mecânica Clássica  (Requisito fraco)`n"

# Step 4: Clear the stray B/C values left behind on row 12 (the row that
# used to hold "Programa resumido:" + a misplaced professor name and is
# now just the "Docentes responsaveis:" label).
$ws.Range("B12:C12").ClearContents()

# Step 5: Row heights - rows 11 and 12 lose the custom height they had
# before the table was restructured; row 17 ("Programa resumido:") gains one.
$ws.Rows(11).AutoFit()
$ws.Rows(12).AutoFit()
$ws.Rows(10).RowHeight = 60
$ws.Rows(17).RowHeight = 60
$ws.Rows(18).RowHeight = 60
$ws.Rows(19).RowHeight = 120
$ws.Rows(20).RowHeight = 120
$ws.Rows(22).RowHeight = 60
$ws.Rows(23).RowHeight = 60
$ws.Rows(24).RowHeight = 60
$ws.Rows(25).RowHeight = 120
$ws.Rows(27).RowHeight = 30

Write-Host "Done"
